$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111749883
$ws.Range("B2").Value = 78107
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("P2").Value = "Älgsjöhåll (Älgsjöhåll), Ög"
$ws.Range("Q2").Value = 575336.5075504743
$ws.Range("R2").Value = 6509789.003789719
$ws.Range("S2").Value = 1

$ws.Range("A3").Value = 111749897
$ws.Range("B3").Value = 78107
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6453
$ws.Range("F3").Value = "Vedskivlav"
$ws.Range("G3").Value = "Hertelidea botryosa"
$ws.Range("H3").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("P3").Value = "Älgsjöhåll (Älgsjöhåll), Ög"
$ws.Range("Q3").Value = 575336.6687912485
$ws.Range("R3").Value = 6509780.695668718
$ws.Range("S3").Value = 1

$ws.Range("A4").Value = 111749006
$ws.Range("B4").Value = 8377
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 106545
$ws.Range("F4").Value = "Mindre märgborre"
$ws.Range("G4").Value = "Tomicus minor"
$ws.Range("H4").Value = "(Hartig, 1834)"
$ws.Range("Q4").Value = 575512.2089522779
$ws.Range("R4").Value = 6509825.662577543

$ws.Range("A5").Value = 111747705
$ws.Range("B5").Value = 93067
$ws.Range("E5").Value = 2810
$ws.Range("F5").Value = "Västlig hakmossa"
$ws.Range("G5").Value = "Rhytidiadelphus loreus"
$ws.Range("H5").Value = "(Hedw.) Warnst."
$ws.Range("Q5").Value = 575459.4222356658
$ws.Range("R5").Value = 6509864.113963567
$ws.Range("S5").Value = 2

$ws.Range("A6").Value = 111749343
$ws.Range("P6").Value = "Lilla gruvan (Lilla gruvan), Ög"
$ws.Range("Q6").Value = 575415.2450877089
$ws.Range("R6").Value = 6509807.674603676

$ws.Range("A7").Value = 111749860
$ws.Range("P7").Value = "Älgsjöhåll (Älgsjöhåll), Ög"
$ws.Range("Q7").Value = 575356.6078101217
$ws.Range("R7").Value = 6509772.251964441

$ws.Range("A8").Value = 111749097
$ws.Range("B8").Value = 93388
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 2180
$ws.Range("F8").Value = "Blåmossa"
$ws.Range("G8").Value = "Leucobryum glaucum"
$ws.Range("H8").Value = "(Hedw.) Ångstr."
$ws.Range("P8").Value = "Lilla gruvan (Lilla gruvan), Ög"
$ws.Range("Q8").Value = 575501.7342092508
$ws.Range("R8").Value = 6509775.591426332
$ws.Range("S8").Value = 3

$ws.Range("A9").Value = 111747186
$ws.Range("P9").Value = "Lilla gruvan (Lilla gruvan), Ög"
$ws.Range("Q9").Value = 575435.6246570286
$ws.Range("R9").Value = 6509856.898648335
$ws.Range("S9").Value = 2
